$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The roster only keeps one name now ("_Ferrero Bonet, Carlos", which sat in
# the old last row, 27). Delete the 23 student rows above it (rows 4-26); the
# surviving row shifts up to become the new row 4, and the sheet's used range
# shrinks from A1:XFC27 down to A1:XFC4.
$ws.Range("A4:A26").EntireRow.Delete() | Out-Null

# Leave the selection where the author's last save left it.
$ws.Range("A13").Select() | Out-Null
